$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Class Statistics summary numbers (K/L columns) ---
$ws.Range("L6").Value = 44
$ws.Range("L7").Value = 13
$ws.Range("L8").Value = 119
$ws.Range("L9").Value = "25.0%"
$ws.Range("L10").Value = "33.1%"

# --- Row 15: PHARMACOLOGY session 2 A1 became Recorded (was Not Recorded) ---
# Copy cell formatting (fill/font) from a known "Recorded" (green) row so the
# style index matches what Excel would assign for the Recorded status band.
$ws.Range("A8:I8").Copy()
$ws.Range("A15:I15").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("G15").Value = "nancy.abdelshafy@med.asu.edu.eg"
$ws.Range("H15").Value = "49/216"
$ws.Range("I15").Value = "Recorded"

# --- Row 15 group-statistics row (O/P/R/S columns) ---
$ws.Range("O15").Value = 7
$ws.Range("P15").Value = 1
$ws.Range("R15").Value = "31.8%"
$ws.Range("S15").Value = "28.8%"

# --- Row 16 group-statistics row (O/Q/R/S columns) ---
$ws.Range("O16").Value = 6
$ws.Range("Q16").Value = 14
$ws.Range("R16").Value = "27.3%"
$ws.Range("S16").Value = "29.2%"

# --- Row 37: PHARMACOLOGY session 2 A2 became Recorded (was Pending) ---
$ws.Range("A8:I8").Copy()
$ws.Range("A37:I37").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("G37").Value = "nancy.abdelshafy@med.asu.edu.eg"
$ws.Range("H37").Value = "2/217"
$ws.Range("I37").Value = "Recorded"

# --- G column (Recorded By) email list re-orderings ---
$ws.Range("G2").Value = "rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G10").Value = "maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
$ws.Range("G18").Value = "yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G19").Value = "naema.gomaa@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G24").Value = "rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G32").Value = "maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
$ws.Range("G40").Value = "yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G41").Value = "naema.gomaa@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G46").Value = "hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
$ws.Range("G54").Value = "basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg"
$ws.Range("G62").Value = "shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G68").Value = "hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
$ws.Range("G76").Value = "basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg"
$ws.Range("G84").Value = "shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G90").Value = "Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G96").Value = "aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G98").Value = "basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg"
$ws.Range("G106").Value = "nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G107").Value = "Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G112").Value = "Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G118").Value = "aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G120").Value = "basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg"
$ws.Range("G128").Value = "nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G129").Value = "Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G134").Value = "hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G142").Value = "basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G150").Value = "naema.gomaa@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G156").Value = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G164").Value = "basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G172").Value = "naema.gomaa@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
